$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 35, shifting existing rows 35-40 down to 36-41.
$ws.Rows.Item(35).Insert()

# New row 35 replicates the same record pattern with updated weekly values.
$ws.Cells.Item(35, 1).Value = 1
$ws.Cells.Item(35, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(35, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(35, 4).Value = 44474
$ws.Cells.Item(35, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(35, 5).Value = 15
$ws.Cells.Item(35, 6).Value = 100112012
$ws.Cells.Item(35, 7).Value = "Espinaca"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 270
$ws.Cells.Item(35, 11).Value = 1000
$ws.Cells.Item(35, 12).Value = 1200
$ws.Cells.Item(35, 13).Value = 1100
$ws.Cells.Item(35, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(35, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(35, 16).Value = 367
$ws.Cells.Item(35, 17).Value = 3
$ws.Cells.Item(35, 18).Value = "Hortaliza"

$wb.Save()
